# Update NATMI Mif-Cxcr2 LR-pair sheet with recomputed TPM-based statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 : ECs -> ECs
$ws.Range("G2").Value  = 6.556445
$ws.Range("H2").Value  = 19.669335
$ws.Range("I2").Value  = 0.2003595613103873
$ws.Range("J2").Value  = 0.2003595613103873
$ws.Range("K2").Value  = 2
$ws.Range("L2").Value  = 0.6666666666666666
$ws.Range("M2").Value  = 0.01650666666666667
$ws.Range("N2").Value  = 0.04952
$ws.Range("O2").Value  = 0.795859985214233
$ws.Range("P2").Value  = 0.795859985214233
$ws.Range("Q2").Value  = 0.1082250521333333
$ws.Range("R2").Value  = 0.9740254692
$ws.Range("S2").Value  = 0.159458157502015
$ws.Range("T2").Value  = 0.159458157502015

# Row 3 : ECs -> FAPs
$ws.Range("G3").Value  = 6.556445
$ws.Range("H3").Value  = 19.669335
$ws.Range("I3").Value  = 0.2003595613103873
$ws.Range("J3").Value  = 0.2003595613103873
$ws.Range("O3").Value  = 0.2041400147857671
$ws.Range("P3").Value  = 0.2041400147857671
$ws.Range("Q3").Value  = 0.02775998813
$ws.Range("R3").Value  = 0.24983989317
$ws.Range("S3").Value  = 0.04090140380837227
$ws.Range("T3").Value  = 0.04090140380837227

# Row 4 : FAPs -> ECs
$ws.Range("I4").Value  = 0.3842320902647997
$ws.Range("J4").Value  = 0.3842320902647997
$ws.Range("K4").Value  = 2
$ws.Range("L4").Value  = 0.6666666666666666
$ws.Range("M4").Value  = 0.01650666666666667
$ws.Range("N4").Value  = 0.04952
$ws.Range("O4").Value  = 0.795859985214233
$ws.Range("P4").Value  = 0.795859985214233
$ws.Range("Q4").Value  = 0.2075445650222222
$ws.Range("R4").Value  = 1.8679010852
$ws.Range("S4").Value  = 0.3057949456769773
$ws.Range("T4").Value  = 0.3057949456769773

# Row 5 : FAPs -> FAPs
$ws.Range("I5").Value  = 0.3842320902647997
$ws.Range("J5").Value  = 0.3842320902647997
$ws.Range("O5").Value  = 0.2041400147857671
$ws.Range("P5").Value  = 0.2041400147857671
$ws.Range("S5").Value  = 0.0784371445878224
$ws.Range("T5").Value  = 0.07843714458782242

# Row 6 : MuSCs -> ECs
$ws.Range("G6").Value  = 13.59357133333334
$ws.Range("I6").Value  = 0.4154083484248129
$ws.Range("J6").Value  = 0.415408348424813
$ws.Range("K6").Value  = 2
$ws.Range("L6").Value  = 0.6666666666666666
$ws.Range("M6").Value  = 0.01650666666666667
$ws.Range("N6").Value  = 0.04952
$ws.Range("O6").Value  = 0.795859985214233
$ws.Range("P6").Value  = 0.795859985214233
$ws.Range("Q6").Value  = 0.2243845508088889
$ws.Range("R6").Value  = 2.01946095728
$ws.Range("S6").Value  = 0.3306068820352406
$ws.Range("T6").Value  = 0.3306068820352406

# Row 7 : MuSCs -> FAPs
$ws.Range("G7").Value  = 13.59357133333334
$ws.Range("I7").Value  = 0.4154083484248129
$ws.Range("J7").Value  = 0.415408348424813
$ws.Range("O7").Value  = 0.2041400147857671
$ws.Range("P7").Value  = 0.2041400147857671
$ws.Range("Q7").Value  = 0.05755518102533334
$ws.Range("R7").Value  = 0.517996629228
$ws.Range("S7").Value  = 0.08480146638957241
$ws.Range("T7").Value  = 0.08480146638957242
